# Update source (col A) and series_code (col K) values on the "timeseries" sheet.
# "umar" -> "umar, surs"; series_code prefix "UMAR--" -> "UMAR-SURS--"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 6; $r++) {
    $sourceCell = $ws.Cells.Item($r, 1)   # column A
    if ($sourceCell.Value2 -eq "umar") {
        $sourceCell.Value2 = "umar, surs"
    }

    $codeCell = $ws.Cells.Item($r, 11)    # column K
    $codeVal = $codeCell.Value2
    if ($codeVal -like "UMAR--*") {
        $codeCell.Value2 = $codeVal -replace "^UMAR--", "UMAR-SURS--"
    }
}

$ws.Range("K7").Select() | Out-Null
